# Edit script applying the documented diff to NC2022-analysis.docx
$d = $word.ActiveDocument

# --- Paragraph 19 (1-indexed): append TODO note about GS to the "measures of bias" paragraph ---
$p18 = $d.Paragraphs(19)
$p18EndNoMark = $p18.Range.End - 1
$ins = $d.Range($p18EndNoMark, $p18EndNoMark)
$ins.InsertAfter("[")
$ins2 = $d.Range($ins.End, $ins.End)
$ins2.InsertAfter("TODO: Discuss GS – over 3 SEM off mean.]")
$ins2.Font.HighlightColorIndex = 7

# --- Paragraph 21 (1-indexed): rewrite the declination paragraph ---
$p20 = $d.Paragraphs(21)
$p20EndNoMark = $p20.Range.End - 1
$p20Target = $d.Range($p20.Range.Start, $p20EndNoMark)
$p20Target.Text = "The difference for declination is also small (0.0798 or less than one tenth of a degree). "

# --- Paragraph 25 (1-indexed): rewrite the responsiveness paragraph, with subscripted "d" ---
$p24 = $d.Paragraphs(25)
$p24EndNoMark = $p24.Range.End - 1
$p24Target = $d.Range($p24.Range.Start, $p24EndNoMark)
$p24Target.Text = "The absolute differences between the measures of responsiveness – big ‘R’, little ‘r’, and the number of responsive districts (Rd) – for the composite and the means for the individual elections are small (0.3774, 0.0212, and 0.0663 for values typically in the low single digits), and the differences are all roughly one standard error or less. "

$p24Full = $d.Paragraphs(25).Range
$findRange = $p24Full.Duplicate
$found = $findRange.Find.Execute("(Rd)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dPos = $findRange.Start + 2
$dRange = $d.Range($dPos, $dPos + 1)
$dRange.Font.Subscript = $true

Write-Output "done"
